$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 3
$ws.Range("F6").Value = 2
$ws.Range("F8").Value = -3
$ws.Range("F12").Value = -3
$ws.Range("F14").Value = 5
